# WIP - Update function insert new record electric
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide / resize column A (RoomID) ---
# (target stored width is 12.6640625; the engine quantizes ColumnWidth to
# steps of 1/6 internally, so 11.833333 is the input that lands on the
# nearest representable stored width, 12.666666666666666)
$colA = $ws.Columns.Item(1)
$colA.Hidden = $false
$colA.ColumnWidth = 11.833333

# --- OldNumber column (B2:B11) gets vertically centered too (in addition to the
#     existing horizontal centering) ---
$ws.Range("B2:B11").HorizontalAlignment = -4108
$ws.Range("B2:B11").VerticalAlignment = -4108

# --- Existing rows (2-11) are marked as processed: Status column G -> 1 ---
$ws.Range("G2:G11").Value = 1

# --- Room 308 (row 11) meter reading got corrected ---
$ws.Cells.Item(11, 4).Value = 21
$ws.Cells.Item(11, 5).Value = 61
$ws.Cells.Item(11, 6).Value = 151000

# --- Insert the new "Bếp" record as row 12 ---
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "Bếp"
$ws.Cells.Item(12, 3).Value = "2019-04"
$ws.Cells.Item(12, 4).Value = 1550
$ws.Cells.Item(12, 5).Value = 1705
$ws.Cells.Item(12, 6).Value = 365000
$ws.Cells.Item(12, 7).Value = 1

# Give the new row the same look as the rest of column A (comma format,
# centered horizontally only - no extra vertical-centering) by copying the
# existing format instead of re-deriving it.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(12, 2).PasteSpecial(-4122)

# --- View updates: selection moves to the newly inserted cell, scroll resets ---
$ws.Range("B12").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
